$wb = $excel.ActiveWorkbook

# The edit: switch focus to the "Repayment Schedule" sheet and insert a new
# (blank) column before column N, shifting the old "Late" / "Outstanding"
# columns one place to the right (N -> O, O -> P, P -> Q).
$ws3 = $wb.Worksheets.Item("Repayment Schedule")
$ws3.Activate()

$ws3.Columns("N:N").Insert()

# The newly inserted column picks up the width of its left neighbour (column M).
$ws3.Columns("N:N").ColumnWidth = 10.3

# Selection ends up on R10 on the Repayment Schedule sheet.
$ws3.Range("R10").Select() | Out-Null
